$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1251056.6
$ws.Range("J9").Value = 812.75
$ws.Range("L9").Value = 812.75
$ws.Range("N9").Value = -1150.75
$ws.Range("H11").Value = 670.7059
$ws.Range("I11").Value = 670.7059
$ws.Range("K11").Value = 670.7059
$ws.Range("M11").Value = -530.7059
$ws.Range("H17").Value = 1093.4865
$ws.Range("J17").Value = 1093.4865
$ws.Range("L17").Value = 3280.4595
$ws.Range("N17").Value = -3616.4595
$ws.Range("H74").Value = 7747.933
$ws.Range("I74").Value = 4121.9
$ws.Range("K74").Value = 4121.9
$ws.Range("M74").Value = -3185.9
$ws.Range("H76").Value = 4923.3335
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685
$ws.Range("H77").Value = 7747.933
$ws.Range("I77").Value = 4121.9
$ws.Range("K77").Value = 20609.5
$ws.Range("M77").Value = -15929.5
$ws.Range("H79").Value = 4923.3335
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908
$ws.Range("H107").Value = 7133.3184
$ws.Range("I107").Value = 7425.381
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 7425.381
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -5505.381
$ws.Range("N107").Value = -4840
$ws.Range("H113").Value = 15991.2
$ws.Range("J113").Value = 4835.3335
$ws.Range("L113").Value = 4835.3335
$ws.Range("N113").Value = -11343.3335
$ws.Range("H129").Value = 2249.5
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H132").Value = 1486.8
$ws.Range("I132").Value = 1057.262
$ws.Range("K132").Value = 3171.786
$ws.Range("M132").Value = -641.7860000000001
$ws.Range("H137").Value = 11440.182
$ws.Range("I137").Value = 13704.667
$ws.Range("J137").Value = 1250
$ws.Range("K137").Value = 41114.001
$ws.Range("L137").Value = 3750
$ws.Range("M137").Value = -38564.001
$ws.Range("N137").Value = -8850
$ws.Range("H138").Value = 5648
$ws.Range("J138").Value = 6106.268
$ws.Range("L138").Value = 18318.804
$ws.Range("N138").Value = -28598.804

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 59998
$ws.Range("J62").Value = 59998
$ws.Range("L62").Value = 59998
$ws.Range("N62").Value = -61246
$ws.Range("H65").Value = 59998
$ws.Range("J65").Value = 59998
$ws.Range("L65").Value = 179994
$ws.Range("N65").Value = -186234
$ws.Range("H74").Value = 3600.5293
$ws.Range("J74").Value = 5247.8887
$ws.Range("L74").Value = 5247.8887
$ws.Range("N74").Value = -6995.8887
$ws.Range("H77").Value = 3600.5293
$ws.Range("J77").Value = 5247.8887
$ws.Range("L77").Value = 26239.4435
$ws.Range("N77").Value = -34975.4435
$ws.Range("H118").Value = 220666.67
$ws.Range("J118").Value = 220666.67
$ws.Range("L118").Value = 220666.67
$ws.Range("N118").Value = -223980.67

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 62503310
$ws.Range("I86").Value = 2999
$ws.Range("J86").Value = 71431930
$ws.Range("K86").Value = 2999
$ws.Range("L86").Value = 71431930
$ws.Range("M86").Value = -1876
$ws.Range("N86").Value = -71434176
$ws.Range("H89").Value = 62503310
$ws.Range("I89").Value = 2999
$ws.Range("J89").Value = 71431930
$ws.Range("K89").Value = 14995
$ws.Range("L89").Value = 357159650
$ws.Range("M89").Value = -9379
$ws.Range("N89").Value = -357170882
$ws.Range("H94").Value = 89546.82000000001
$ws.Range("I94").Value = 112498.09
$ws.Range("K94").Value = 112498.09
$ws.Range("M94").Value = -112047.09
$ws.Range("H105").Value = 2282.261
$ws.Range("I105").Value = 1817.75
$ws.Range("K105").Value = 1817.75
$ws.Range("M105").Value = -70.75
$ws.Range("H134").Value = 2220.5854
$ws.Range("I134").Value = 2168.9143
$ws.Range("J134").Value = 2522
$ws.Range("K134").Value = 6506.742899999999
$ws.Range("L134").Value = 7566
$ws.Range("M134").Value = -3971.742899999999
$ws.Range("N134").Value = -12636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1555.8
$ws.Range("I16").Value = 694.75
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 694.75
$ws.Range("L16").Value = 5000
$ws.Range("M16").Value = -407.75
$ws.Range("N16").Value = -5574
$ws.Range("H31").Value = 3175.8215
$ws.Range("I31").Value = 1242.2632
$ws.Range("J31").Value = 7257.778
$ws.Range("K31").Value = 1242.2632
$ws.Range("L31").Value = 7257.778
$ws.Range("M31").Value = -947.2632000000001
$ws.Range("N31").Value = -7847.778
$ws.Range("H34").Value = 3175.8215
$ws.Range("I34").Value = 1242.2632
$ws.Range("J34").Value = 7257.778
$ws.Range("K34").Value = 1242.2632
$ws.Range("L34").Value = 7257.778
$ws.Range("M34").Value = -1040.2632
$ws.Range("N34").Value = -7661.778
$ws.Range("H58").Value = 2741.1428
$ws.Range("J58").Value = 3757
$ws.Range("L58").Value = 3757
$ws.Range("N58").Value = -4163
$ws.Range("H113").Value = 1555.8
$ws.Range("I113").Value = 694.75
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 694.75
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 1475.25
$ws.Range("N113").Value = -9340
$ws.Range("H136").Value = 2741.1428
$ws.Range("J136").Value = 3757
$ws.Range("L136").Value = 11271
$ws.Range("N136").Value = -16371

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 800
$ws.Range("I63").Value = 800
$ws.Range("K63").Value = 2400
$ws.Range("M63").Value = -1651
$ws.Range("H66").Value = 800
$ws.Range("I66").Value = 800
$ws.Range("K66").Value = 7200
$ws.Range("M66").Value = -3456
$ws.Range("H113").Value = 1311
$ws.Range("I113").Value = 537.5
$ws.Range("J113").Value = 1929.8
$ws.Range("K113").Value = 1612.5
$ws.Range("L113").Value = 5789.4
$ws.Range("M113").Value = 557.5
$ws.Range("N113").Value = -10129.4
$ws.Range("H137").Value = 1676.6666
$ws.Range("I137").Value = 1676.6666
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5029.9998
$ws.Range("L137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = 70.0002000000004
$ws.Range("H139").Value = 3619.5186
$ws.Range("I139").Value = 1430.6364
$ws.Range("K139").Value = 4291.9092
$ws.Range("M139").Value = 848.0907999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9294.151
$ws.Range("I70").Value = 9813.571
$ws.Range("J70").Value = 9154.308000000001
$ws.Range("K70").Value = 9813.571
$ws.Range("L70").Value = 9154.308000000001
$ws.Range("M70").Value = -9543.571
$ws.Range("N70").Value = -9694.308000000001
$ws.Range("H73").Value = 9294.151
$ws.Range("I73").Value = 9813.571
$ws.Range("J73").Value = 9154.308000000001
$ws.Range("K73").Value = 9813.571
$ws.Range("L73").Value = 9154.308000000001
$ws.Range("M73").Value = -8877.571
$ws.Range("N73").Value = -11026.308
$ws.Range("H80").Value = 12345.786
$ws.Range("I80").Value = 17934.428
$ws.Range("J80").Value = 6757.143
$ws.Range("K80").Value = 17934.428
$ws.Range("L80").Value = 6757.143
$ws.Range("M80").Value = -16936.428
$ws.Range("N80").Value = -8753.143
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").ClearContents()
$ws.Range("N82").Value = 0
$ws.Range("H83").Value = 12345.786
$ws.Range("I83").Value = 17934.428
$ws.Range("J83").Value = 6757.143
$ws.Range("K83").Value = 89672.14
$ws.Range("L83").Value = 33785.715
$ws.Range("M83").Value = -84680.14
$ws.Range("N83").Value = -43769.715
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").ClearContents()
$ws.Range("N85").Value = 0
$ws.Range("H113").Value = 28577548
$ws.Range("I113").Value = 62501668
$ws.Range("K113").Value = 62501668
$ws.Range("M113").Value = -62499498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 305.45456
$ws.Range("I55").Value = 316
$ws.Range("J55").Value = 287
$ws.Range("K55").Value = 316
$ws.Range("L55").Value = 287
$ws.Range("M55").Value = -143
$ws.Range("N55").Value = -633
$ws.Range("H61").Value = 4914.84
$ws.Range("I61").Value = 1780.5
$ws.Range("J61").Value = 7808.077
$ws.Range("K61").Value = 1780.5
$ws.Range("L61").Value = 7808.077
$ws.Range("M61").Value = -1578.5
$ws.Range("N61").Value = -8212.077000000001
$ws.Range("H93").Value = 3198.72
$ws.Range("J93").Value = 5386.091
$ws.Range("L93").Value = 5386.091
$ws.Range("N93").Value = -7882.091
$ws.Range("H113").Value = 4914.84
$ws.Range("I113").Value = 1780.5
$ws.Range("J113").Value = 7808.077
$ws.Range("K113").Value = 1780.5
$ws.Range("L113").Value = 7808.077
$ws.Range("M113").Value = 389.5
$ws.Range("N113").Value = -12148.077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 170566
$ws.Range("J120").Value = 170566
$ws.Range("L120").Value = 170566
$ws.Range("N120").Value = -180242
$ws.Range("H127").Value = 89995
$ws.Range("J127").Value = 89995
$ws.Range("L127").Value = 89995
$ws.Range("N127").Value = -99915
$ws.Range("H135").Value = 57857.5
$ws.Range("J135").Value = 57857.5
$ws.Range("L135").Value = 57857.5
